$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7602039999999999
$ws.Range("H2").Value = 2.280612
$ws.Range("I2").Value = 0.9081302434927829
$ws.Range("J2").Value = 0.9140443705567521
$ws.Range("M2").Value = 55.96137100000001
$ws.Range("N2").Value = 167.884113
$ws.Range("O2").Value = 0.7743971326715885
$ws.Range("P2").Value = 0.7791282367296697
$ws.Range("Q2").Value = 42.542058079684
$ws.Range("R2").Value = 382.878522717156
$ws.Range("S2").Value = 0.7032534566531626
$ws.Range("T2").Value = 0.7121577787245631
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7602039999999999
$ws.Range("H3").Value = 2.280612
$ws.Range("I3").Value = 0.9081302434927829
$ws.Range("J3").Value = 0.9140443705567521
$ws.Range("O3").Value = 0.2047066592056736
$ws.Range("P3").Value = 0.20595729465514
$ws.Range("Q3").Value = 11.245706134244
$ws.Range("R3").Value = 101.211355208196
$ws.Range("S3").Value = 0.1859003082690425
$ws.Range("T3").Value = 0.1882541057546289
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7602039999999999
$ws.Range("H4").Value = 2.280612
$ws.Range("I4").Value = 0.9081302434927829
$ws.Range("J4").Value = 0.9140443705567521
$ws.Range("M4").Value = 0.07806133333333333
$ws.Range("N4").Value = 0.234184
$ws.Range("O4").Value = 0.001080217865031477
$ws.Range("P4").Value = 0.001086817351146865
$ws.Range("Q4").Value = 0.05934253784533332
$ws.Range("R4").Value = 0.5340828406079999
$ws.Range("S4").Value = 0.0009809785127962893
$ws.Range("T4").Value = 0.0009933992816391929
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7602039999999999
$ws.Range("H5").Value = 2.280612
$ws.Range("I5").Value = 0.9081302434927829
$ws.Range("J5").Value = 0.9140443705567521
$ws.Range("M5").Value = 1.316435
$ws.Range("N5").Value = 2.63287
$ws.Range("O5").Value = 0.01821691411649872
$ws.Range("P5").Value = 0.01221880572248337
$ws.Range("Q5").Value = 1.00075915274
$ws.Range("R5").Value = 6.004554916439999
$ws.Range("S5").Value = 0.0165433306523031
$ws.Range("T5").Value = 0.01116853058556256
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7602039999999999
$ws.Range("H6").Value = 2.280612
$ws.Range("I6").Value = 0.9081302434927829
$ws.Range("J6").Value = 0.9140443705567521
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1155563333333333
$ws.Range("N6").Value = 0.346669
$ws.Range("O6").Value = 0.001599076141207756
$ws.Range("P6").Value = 0.001608845541560194
$ws.Range("Q6").Value = 0.08784638682533331
$ws.Range("R6").Value = 0.7906174814279999
$ws.Range("S6").Value = 0.001452169405478499
$ws.Range("T6").Value = 0.001470556210358425
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.016249
$ws.Range("H7").Value = 0.032498
$ws.Range("I7").Value = 0.01941085330584189
$ws.Range("J7").Value = 0.01302484331151171
$ws.Range("M7").Value = 55.96137100000001
$ws.Range("N7").Value = 167.884113
$ws.Range("O7").Value = 0.7743971326715885
$ws.Range("P7").Value = 0.7791282367296697
$ws.Range("Q7").Value = 0.9093163173790001
$ws.Range("R7").Value = 5.455897904274
$ws.Range("S7").Value = 0.01503170914275279
$ws.Range("T7").Value = 0.01014802320297835
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.016249
$ws.Range("H8").Value = 0.032498
$ws.Range("I8").Value = 0.01941085330584189
$ws.Range("J8").Value = 0.01302484331151171
$ws.Range("O8").Value = 0.2047066592056736
$ws.Range("P8").Value = 0.20595729465514
$ws.Range("Q8").Value = 0.240371635739
$ws.Range("R8").Value = 1.442229814434
$ws.Range("S8").Value = 0.003973530932570299
$ws.Range("T8").Value = 0.002682561491746046
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.016249
$ws.Range("H9").Value = 0.032498
$ws.Range("I9").Value = 0.01941085330584189
$ws.Range("J9").Value = 0.01302484331151171
$ws.Range("M9").Value = 0.07806133333333333
$ws.Range("N9").Value = 0.234184
$ws.Range("O9").Value = 0.001080217865031477
$ws.Range("P9").Value = 0.001086817351146865
$ws.Range("Q9").Value = 0.001268418605333333
$ws.Range("R9").Value = 0.007610511632
$ws.Range("S9").Value = 0.00002096795051647572
$ws.Range("T9").Value = 0.00001415562570692011
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.016249
$ws.Range("H10").Value = 0.032498
$ws.Range("I10").Value = 0.01941085330584189
$ws.Range("J10").Value = 0.01302484331151171
$ws.Range("M10").Value = 1.316435
$ws.Range("N10").Value = 2.63287
$ws.Range("O10").Value = 0.01821691411649872
$ws.Range("P10").Value = 0.01221880572248337
$ws.Range("Q10").Value = 0.021390752315
$ws.Range("R10").Value = 0.08556300926
$ws.Range("S10").Value = 0.0003536058476004771
$ws.Range("T10").Value = 0.0001591480299891485
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.016249
$ws.Range("H11").Value = 0.032498
$ws.Range("I11").Value = 0.01941085330584189
$ws.Range("J11").Value = 0.01302484331151171
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1155563333333333
$ws.Range("N11").Value = 0.346669
$ws.Range("O11").Value = 0.001599076141207756
$ws.Range("P11").Value = 0.001608845541560194
$ws.Range("Q11").Value = 0.001877674860333333
$ws.Range("R11").Value = 0.011266049162
$ws.Range("S11").Value = 0.00003103943240185546
$ws.Range("T11").Value = 0.00002095496109124572
$ws.Range("G12").Value = 0.06065600000000001
$ws.Range("H12").Value = 0.181968
$ws.Range("I12").Value = 0.07245890320137523
$ws.Range("J12").Value = 0.07293078613173619
$ws.Range("M12").Value = 55.96137100000001
$ws.Range("N12").Value = 167.884113
$ws.Range("O12").Value = 0.7743971326715885
$ws.Range("P12").Value = 0.7791282367296697
$ws.Range("Q12").Value = 3.394392919376001
$ws.Range("R12").Value = 30.54953627438401
$ws.Range("S12").Value = 0.05611196687567316
$ws.Range("T12").Value = 0.05682243480212826
$ws.Range("G13").Value = 0.06065600000000001
$ws.Range("H13").Value = 0.181968
$ws.Range("I13").Value = 0.07245890320137523
$ws.Range("J13").Value = 0.07293078613173619
$ws.Range("O13").Value = 0.2047066592056736
$ws.Range("P13").Value = 0.20595729465514
$ws.Range("Q13").Value = 0.8972848752160001
$ws.Range("R13").Value = 8.075563876944001
$ws.Range("S13").Value = 0.01483282000406081
$ws.Range("T13").Value = 0.01502062740876499
$ws.Range("G14").Value = 0.06065600000000001
$ws.Range("H14").Value = 0.181968
$ws.Range("I14").Value = 0.07245890320137523
$ws.Range("J14").Value = 0.07293078613173619
$ws.Range("M14").Value = 0.07806133333333333
$ws.Range("N14").Value = 0.234184
$ws.Range("O14").Value = 0.001080217865031477
$ws.Range("P14").Value = 0.001086817351146865
$ws.Range("Q14").Value = 0.004734888234666667
$ws.Range("R14").Value = 0.042613994112
$ws.Range("S14").Value = 0.000078271401718712
$ws.Range("T14").Value = 0.00007926244380075204
$ws.Range("G15").Value = 0.06065600000000001
$ws.Range("H15").Value = 0.181968
$ws.Range("I15").Value = 0.07245890320137523
$ws.Range("J15").Value = 0.07293078613173619
$ws.Range("M15").Value = 1.316435
$ws.Range("N15").Value = 2.63287
$ws.Range("O15").Value = 0.01821691411649872
$ws.Range("P15").Value = 0.01221880572248337
$ws.Range("Q15").Value = 0.07984968136000001
$ws.Range("R15").Value = 0.47909808816
$ws.Range("S15").Value = 0.001319977616595147
$ws.Range("T15").Value = 0.0008911271069316692
$ws.Range("G16").Value = 0.06065600000000001
$ws.Range("H16").Value = 0.181968
$ws.Range("I16").Value = 0.07245890320137523
$ws.Range("J16").Value = 0.07293078613173619
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1155563333333333
$ws.Range("N16").Value = 0.346669
$ws.Range("O16").Value = 0.001599076141207756
$ws.Range("P16").Value = 0.001608845541560194
$ws.Range("Q16").Value = 0.007009184954666668
$ws.Range("R16").Value = 0.063082664592
$ws.Range("S16").Value = 0.0001158673033274014
$ws.Range("T16").Value = 0.0001173343701105238
